$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column L (shifts old L:R -> N:T)
$ws.Range("L1:M1").EntireColumn.Insert()

# Give the new header cells the same formatting as the neighboring header cell (K1)
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Match the new columns' width to column K's width
$ws.Range("L:M").ColumnWidth = $ws.Range("K:K").ColumnWidth

# Set the new header cell values
$ws.Range("L1").Value = "บัตร"
$ws.Range("M1").Value = "IPD/OPD"

# The data validation list (now living at N:O) needs its source formula updated to
# point at the shifted helper column (T instead of R); Insert() shifts the sqref
# automatically but not the formula, so update it in place.
$dv = $ws.Range("N2:O1048576").Validation
$dv.Modify([Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList, [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop, [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween, "=`$T`$2:`$T`$3")
